# Updates the crypto price/volume table to the latest scraped values.
# D-column prices are stored as TEXT (e.g. "58.135.04", "0.150") so they
# must be forced to text (NumberFormat "@") before assignment, otherwise
# Excel auto-coerces numeric-looking strings into numbers and silently
# drops significant trailing zeros / thousand-dot formatting. The style
# is reset back to "Normal" right after so no visible formatting changes
# leak into the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue "D2" '58.135.04'
$ws.Range("E2").Value = '  -0.64%  '
Set-TextValue "D3" '2.286.29'
$ws.Range("E3").Value = '  +0.43%  '
$ws.Range("E4").Value = '  +0.05%  '
Set-TextValue "D5" '534.99'
$ws.Range("E5").Value = '  -1.74%  '
Set-TextValue "D6" '131.39'
$ws.Range("E6").Value = '  +0.69%  '
Set-TextValue "D7" '0.999'
$ws.Range("E7").Value = '  +0.03%  '
Set-TextValue "D8" '0.588'
Set-TextValue "D9" '2.286.50'
$ws.Range("E9").Value = '  +0.44%  '
Set-TextValue "D10" '0.0999'
$ws.Range("E10").Value = '  -1.35%  '
$ws.Range("E11").Value = '  -0.05%  '
Set-TextValue "D12" '0.150'
$ws.Range("E12").Value = '  +0.80%  '
$ws.Range("E13").Value = '  -0.38%  '
Set-TextValue "D14" '23.50'
$ws.Range("E14").Value = '  -0.51%  '
Set-TextValue "D15" '2.691.58'
$ws.Range("E15").Value = '  +0.48%  '
Set-TextValue "D16" '58.018.79'
$ws.Range("E16").Value = '  -0.74%  '
$ws.Range("E17").Value = '  -0.40%  '
Set-TextValue "D18" '2.284.80'
$ws.Range("E18").Value = '  +0.61%  '
Set-TextValue "D19" '10.51'
$ws.Range("E19").Value = '  -1.46%  '
Set-TextValue "D20" '4.19'
$ws.Range("E20").Value = '  -2.25%  '
Set-TextValue "D21" '313.09'
$ws.Range("E21").Value = '  -0.18%  '
Set-TextValue "D22" '6.46'
$ws.Range("E22").Value = '  +0.37%  '
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("E24").Value = '  +0.46%  '
$ws.Range("E25").Value = '  -0.69%  '
Set-TextValue "D26" '1.03'
$ws.Range("E26").Value = '  +2.63%  '
Set-TextValue "D27" '7.97'
$ws.Range("E27").Value = '  -1.71%  '
Set-TextValue "D29" '170.68'
$ws.Range("E29").Value = '  -0.24%  '
$ws.Range("E30").Value = '  -2.11%  '
Set-TextValue "D31" '0.0₃0722'
$ws.Range("E31").Value = '  +0.54%  '
$ws.Range("B32").Value = 'SuiNetwork'
$ws.Range("C32").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue "D32" '1.08'
$ws.Range("E32").Value = '  +0.04%  '
$ws.Range("B33").Value = 'Aptos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue "D33" '5.76'
$ws.Range("E33").Value = '  -0.06%  '
Set-TextValue "D34" '0.380'
$ws.Range("E34").Value = '  -0.26%  '
Set-TextValue "D36" '17.84'
$ws.Range("E36").Value = '  +0.71%  '
$ws.Range("E37").Value = '  -0.05%  '
Set-TextValue "D38" '1.24'
$ws.Range("E38").Value = '  -1.20%  '
Set-TextValue "D39" '3.90'
$ws.Range("E39").Value = '  -1.18%  '
$ws.Range("E40").Value = '  -1.28%  '
Set-TextValue "D41" '139.91'
$ws.Range("E41").Value = '  -0.31%  '
Set-TextValue "D42" '287.41'
$ws.Range("E42").Value = '  -4.07%  '
Set-TextValue "D43" '3.44'
$ws.Range("E43").Value = '  -0.14%  '
$ws.Range("E44").Value = '  +0.84%  '
Set-TextValue "D46" '0.553'
$ws.Range("E46").Value = '  +0.77%  '
Set-TextValue "D47" '18.10'
$ws.Range("E47").Value = '  -1.01%  '
$ws.Range("E48").Value = '  -1.33%  '
Set-TextValue "D49" '10.95'
$ws.Range("E49").Value = '  -0.50%  '
$ws.Range("E50").Value = '  -0.21%  '
$ws.Range("E51").Value = '  +1.43%  '
